$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D:G to be treated as text so numeric-looking strings
# (prices, percentages, dates, hours) are preserved exactly as typed,
# matching the original inlineStr cell contents.
$ws.Range("D2:G51").NumberFormat = "@"

$ws.Range("D2").Value = '290.87'
$ws.Range("E2").Value = '-6.99%'
$ws.Range("F2").Value = '14-2-2023'
$ws.Range("G2").Value = '1'

$ws.Range("D3").Value = '40.91'
$ws.Range("E3").Value = '1.83%'
$ws.Range("F3").Value = '14-2-2023'
$ws.Range("G3").Value = '1'

$ws.Range("D4").Value = '5.034'
$ws.Range("E4").Value = '-2.25%'
$ws.Range("F4").Value = '14-2-2023'
$ws.Range("G4").Value = '1'

$ws.Range("D5").Value = '0.07354'
$ws.Range("E5").Value = '-3.06%'
$ws.Range("F5").Value = '14-2-2023'
$ws.Range("G5").Value = '1'

$ws.Range("D6").Value = '4.296'
$ws.Range("E6").Value = '-0.99%'
$ws.Range("F6").Value = '14-2-2023'
$ws.Range("G6").Value = '1'

$ws.Range("D7").Value = '1.562'
$ws.Range("E7").Value = '-8.39%'
$ws.Range("F7").Value = '14-2-2023'
$ws.Range("G7").Value = '1'

$ws.Range("D8").Value = '0.9262'
$ws.Range("E8").Value = '-1.15%'
$ws.Range("F8").Value = '14-2-2023'
$ws.Range("G8").Value = '1'

$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").Value = '0.1179'
$ws.Range("E9").Value = '-5.98%'
$ws.Range("F9").Value = '14-2-2023'
$ws.Range("G9").Value = '1'

$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = '0.1741'
$ws.Range("E10").Value = '-5.11%'
$ws.Range("F10").Value = '14-2-2023'
$ws.Range("G10").Value = '1'

$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = '0.08621'
$ws.Range("E11").Value = '-4.75%'
$ws.Range("F11").Value = '14-2-2023'
$ws.Range("G11").Value = '1'

$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = '0.04169'
$ws.Range("E12").Value = '2.02%'
$ws.Range("F12").Value = '14-2-2023'
$ws.Range("G12").Value = '1'

$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = '0.1053'
$ws.Range("E13").Value = '-0.22%'
$ws.Range("F13").Value = '14-2-2023'
$ws.Range("G13").Value = '1'

$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = '0.001275'
$ws.Range("E14").Value = '-0.98%'
$ws.Range("F14").Value = '14-2-2023'
$ws.Range("G14").Value = '1'

$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").Value = '0.006043'
$ws.Range("E15").Value = '2.45%'
$ws.Range("F15").Value = '14-2-2023'
$ws.Range("G15").Value = '1'

$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").Value = '3.376'
$ws.Range("E16").Value = '0.58%'
$ws.Range("F16").Value = '14-2-2023'
$ws.Range("G16").Value = '1'

$ws.Range("B17").Value = 'BTSEToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D17").Value = '2.378'
$ws.Range("E17").Value = '-1.94%'
$ws.Range("F17").Value = '14-2-2023'
$ws.Range("G17").Value = '1'

$ws.Range("D18").Value = '0.3288'
$ws.Range("E18").Value = '-2.06%'
$ws.Range("F18").Value = '14-2-2023'
$ws.Range("G18").Value = '1'

$ws.Range("D19").Value = '7.696'
$ws.Range("E19").Value = '-9.44%'
$ws.Range("F19").Value = '14-2-2023'
$ws.Range("G19").Value = '1'

$ws.Range("D20").Value = '0.1381'
$ws.Range("E20").Value = '2.77%'
$ws.Range("F20").Value = '14-2-2023'
$ws.Range("G20").Value = '1'

$ws.Range("D21").Value = '0.2883'
$ws.Range("E21").Value = '5.67%'
$ws.Range("F21").Value = '14-2-2023'
$ws.Range("G21").Value = '1'

$ws.Range("D22").Value = '0.03898'
$ws.Range("E22").Value = '-3.37%'
$ws.Range("F22").Value = '14-2-2023'
$ws.Range("G22").Value = '1'

$ws.Range("D23").Value = '0.001260'
$ws.Range("E23").Value = '-0.51%'
$ws.Range("F23").Value = '14-2-2023'
$ws.Range("G23").Value = '1'

$ws.Range("D24").Value = '0.003802'
$ws.Range("E24").Value = '-6.33%'
$ws.Range("F24").Value = '14-2-2023'
$ws.Range("G24").Value = '1'

$ws.Range("D25").Value = '0.0001281'
$ws.Range("E25").Value = '0.58%'
$ws.Range("F25").Value = '14-2-2023'
$ws.Range("G25").Value = '1'

$ws.Range("D26").Value = '0.0003724'
$ws.Range("E26").Value = '-95.05%'
$ws.Range("F26").Value = '14-2-2023'
$ws.Range("G26").Value = '1'

$ws.Range("F27").Value = '14-2-2023'
$ws.Range("G27").Value = '1'

$ws.Range("F28").Value = '14-2-2023'
$ws.Range("G28").Value = '1'

$ws.Range("F29").Value = '14-2-2023'
$ws.Range("G29").Value = '1'

$ws.Range("F30").Value = '14-2-2023'
$ws.Range("G30").Value = '1'

$ws.Range("F31").Value = '14-2-2023'
$ws.Range("G31").Value = '1'

$ws.Range("F32").Value = '14-2-2023'
$ws.Range("G32").Value = '1'

$ws.Range("F33").Value = '14-2-2023'
$ws.Range("G33").Value = '1'

$ws.Range("F34").Value = '14-2-2023'
$ws.Range("G34").Value = '1'

$ws.Range("F35").Value = '14-2-2023'
$ws.Range("G35").Value = '1'

$ws.Range("F36").Value = '14-2-2023'
$ws.Range("G36").Value = '1'

$ws.Range("F37").Value = '14-2-2023'
$ws.Range("G37").Value = '1'

$ws.Range("D38").Value = '0.02323'
$ws.Range("E38").Value = '-6.06%'
$ws.Range("F38").Value = '14-2-2023'
$ws.Range("G38").Value = '1'

$ws.Range("D39").Value = '0.05006'
$ws.Range("E39").Value = '-3.79%'
$ws.Range("F39").Value = '14-2-2023'
$ws.Range("G39").Value = '1'

$ws.Range("D40").Value = '0.005894'
$ws.Range("E40").Value = '177.50%'
$ws.Range("F40").Value = '14-2-2023'
$ws.Range("G40").Value = '1'

$ws.Range("D41").Value = '0.007686'
$ws.Range("E41").Value = '-1.19%'
$ws.Range("F41").Value = '14-2-2023'
$ws.Range("G41").Value = '1'

$ws.Range("D42").Value = '0.1285'
$ws.Range("E42").Value = '-1.04%'
$ws.Range("F42").Value = '14-2-2023'
$ws.Range("G42").Value = '1'

$ws.Range("D43").Value = '0.007348'
$ws.Range("E43").Value = '-4.52%'
$ws.Range("F43").Value = '14-2-2023'
$ws.Range("G43").Value = '1'

$ws.Range("D44").Value = '0.007092'
$ws.Range("E44").Value = '-12.48%'
$ws.Range("F44").Value = '14-2-2023'
$ws.Range("G44").Value = '1'

$ws.Range("D45").Value = '0.3138'
$ws.Range("E45").Value = '0.28%'
$ws.Range("F45").Value = '14-2-2023'
$ws.Range("G45").Value = '1'

$ws.Range("D46").Value = '0.00006367'
$ws.Range("E46").Value = '-4.16%'
$ws.Range("F46").Value = '14-2-2023'
$ws.Range("G46").Value = '1'

$ws.Range("D47").Value = '0.00000000751'
$ws.Range("E47").Value = '-0.21%'
$ws.Range("F47").Value = '14-2-2023'
$ws.Range("G47").Value = '1'

$ws.Range("D48").Value = '0.01779'
$ws.Range("E48").Value = '-57.14%'
$ws.Range("F48").Value = '14-2-2023'
$ws.Range("G48").Value = '1'

$ws.Range("D49").Value = '0.00002101'
$ws.Range("E49").Value = '-0.21%'
$ws.Range("F49").Value = '14-2-2023'
$ws.Range("G49").Value = '1'

$ws.Range("D50").Value = '0.0002001'
$ws.Range("E50").Value = '-0.21%'
$ws.Range("F50").Value = '14-2-2023'
$ws.Range("G50").Value = '1'

$ws.Range("F51").Value = '14-2-2023'
$ws.Range("G51").Value = '1'

# Restore the default "Normal" style on the affected range so the
# text number-format tweak above does not leave a stray style index
# on the cells (they had no explicit style before the edit either).
$ws.Range("D2:G51").Style = "Normal"
